$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Row 2 (header2): swap the Yes/No order for the Ecology and Social sub-columns
$tbl.Cell(2, 2).Range.Text = "Yes"
$tbl.Cell(2, 3).Range.Text = "No"
$tbl.Cell(2, 4).Range.Text = "Yes"
$tbl.Cell(2, 5).Range.Text = "No"

# Row 3 (body1): keep the counts aligned with their (now swapped) Yes/No label
$tbl.Cell(3, 2).Range.Text = "3"
$tbl.Cell(3, 3).Range.Text = "26"
$tbl.Cell(3, 4).Range.Text = "4"
$tbl.Cell(3, 5).Range.Text = "9"
